# Rebuild the "Estado de Cuenta" worker table (rows 16-66).
# The previous account-statement rows are replaced by a new, re-sorted data set:
# grouped by worker (in a new order) with periods listed from the most recent
# (2406) down to the oldest (2312), and CANDELARIA TERRAZA LEDESMA's Salario
# Basico corrected from 1300000/908526-style values to 909000.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, TipoDoc(B), NDoc(C), NombreTrabajador(D), PeriodoMora(E), ValorMora(F), SalarioBasico(G)
$data = @(
  @(16, "CC", "45541756", "CANDELARIA TERRAZA LEDESMA", "2401", 46400, 909000),
  @(17, "CC", "45541756", "CANDELARIA TERRAZA LEDESMA", "2312", 46400, 909000),
  @(18, "CC", "7918528", "HARLOS MORALES REVOLLO", "2406", 18560, 1160000),
  @(19, "CC", "7918528", "HARLOS MORALES REVOLLO", "2405", 46400, 1160000),
  @(20, "CC", "7918528", "HARLOS MORALES REVOLLO", "2404", 46400, 1160000),
  @(21, "CC", "7918528", "HARLOS MORALES REVOLLO", "2403", 46400, 1160000),
  @(22, "CC", "7918528", "HARLOS MORALES REVOLLO", "2402", 46400, 1160000),
  @(23, "CC", "7918528", "HARLOS MORALES REVOLLO", "2401", 46400, 1160000),
  @(24, "CC", "7918528", "HARLOS MORALES REVOLLO", "2312", 46400, 1160000),
  @(25, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2406", 18560, 1160000),
  @(26, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2405", 46400, 1160000),
  @(27, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2404", 46400, 1160000),
  @(28, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2403", 46400, 1160000),
  @(29, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2402", 46400, 1160000),
  @(30, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2401", 46400, 1160000),
  @(31, "CC", "1038123871", "ANDRES FELIPE CONTRERAS LOBOS", "2312", 46400, 1160000),
  @(32, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2406", 18560, 908526),
  @(33, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2405", 46400, 908526),
  @(34, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2404", 46400, 908526),
  @(35, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2403", 46400, 908526),
  @(36, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2402", 46400, 908526),
  @(37, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2401", 46400, 908526),
  @(38, "CC", "73154575", "JUAN CARLOS MONTALVO MORALES", "2312", 46400, 908526),
  @(39, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2406", 18560, 1160000),
  @(40, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2405", 46400, 1160000),
  @(41, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2404", 46400, 1160000),
  @(42, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2403", 46400, 1160000),
  @(43, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2402", 46400, 1160000),
  @(44, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2401", 46400, 1160000),
  @(45, "CC", "1047482188", "JULIETH PAOLA MUÑOZ BLANCO", "2312", 46400, 1160000),
  @(46, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2406", 18560, 1160000),
  @(47, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2405", 46400, 1160000),
  @(48, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2404", 46400, 1160000),
  @(49, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2403", 46400, 1160000),
  @(50, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2402", 46400, 1160000),
  @(51, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2401", 46400, 1160000),
  @(52, "CC", "32938821", "SANDRA MILENA CASTRO DE VOZ", "2312", 46400, 1160000),
  @(53, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2406", 18560, 1160000),
  @(54, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2405", 46400, 1160000),
  @(55, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2404", 46400, 1160000),
  @(56, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2403", 46400, 1160000),
  @(57, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2402", 46400, 1160000),
  @(58, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2401", 46400, 1160000),
  @(59, "CC", "1051736916", "YADIRA DEL CARMEN MENCO CASTRO", "2312", 46400, 1160000),
  @(60, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2406", 18560, 1160000),
  @(61, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2405", 46400, 1160000),
  @(62, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2404", 46400, 1160000),
  @(63, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2403", 46400, 1160000),
  @(64, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2402", 46400, 1160000),
  @(65, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2401", 46400, 1160000),
  @(66, "CC", "73102272", "CARLOS MARIANO CASSERES ROMAN", "2312", 46400, 1160000)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]   # B: Tipo Doc Trabajador
  $ws.Cells.Item($r, 3).Value = $row[2]   # C: N Doc Trabajador
  $ws.Cells.Item($r, 4).Value = $row[3]   # D: Nombre Trabajador
  $ws.Cells.Item($r, 5).Value = $row[4]   # E: Periodo Mora
  $ws.Cells.Item($r, 6).Value = $row[5]   # F: Valor Mora
  $ws.Cells.Item($r, 7).Value = $row[6]   # G: Salario Basico
}
